$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.538.99"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").Value = "2.225.25"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'270.17"
$ws.Range("E5").Value = "  +3.82%  "

$ws.Range("D6").Value = "'92.57"
$ws.Range("E6").Value = "  +12.38%  "

$ws.Range("E7").Value = "  -0.97%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +2.34%  "

$ws.Range("D10").Value = "'45.94"
$ws.Range("E10").Value = "  +5.11%  "

$ws.Range("D11").Value = "'0.0931"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").Value = "'8.18"
$ws.Range("E12").Value = "  +16.19%  "

$ws.Range("D13").Value = "'0.105"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("D14").Value = "2.559.85"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").Value = "'15.04"
$ws.Range("E15").Value = "  +3.55%  "

$ws.Range("D16").Value = "2.226.58"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("E17").Value = "  +2.22%  "

$ws.Range("D18").Value = "43.525.66"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("D20").Value = "'6.00"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").Value = "'70.37"
$ws.Range("E21").Value = "  -1.27%  "

$ws.Range("D22").Value = "'2.35"
$ws.Range("E22").Value = "  -0.89%  "

$ws.Range("D23").Value = "'232.43"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").Value = "'9.08"
$ws.Range("E24").Value = "  -2.27%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'11.35"
$ws.Range("E26").Value = "  +5.43%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.50"
$ws.Range("E27").Value = "  +11.10%  "

$ws.Range("E28").Value = "  +5.31%  "

$ws.Range("D29").Value = "'41.07"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("E30").Value = "  +1.96%  "

$ws.Range("D31").Value = "'172.59"
$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").Value = "'0.0922"
$ws.Range("E32").Value = "  +5.55%  "

$ws.Range("D33").Value = "'20.82"
$ws.Range("E33").Value = "  +1.00%  "

$ws.Range("D34").Value = "'5.46"
$ws.Range("E34").Value = "  +2.38%  "

$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("D36").Value = "'0.112"
$ws.Range("E36").Value = "  -3.52%  "

$ws.Range("E37").Value = "  -2.59%  "

$ws.Range("D38").Value = "'4.30"
$ws.Range("E38").Value = "  -4.63%  "

$ws.Range("D39").Value = "'3.59"
$ws.Range("E39").Value = "  +21.72%  "

$ws.Range("D40").Value = "'12.54"
$ws.Range("E40").Value = "  -7.26%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.219"
$ws.Range("E41").Value = "  +8.84%  "

$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").Value = "'2.16"
$ws.Range("E42").Value = "  +2.17%  "

$ws.Range("D43").Value = "'63.21"
$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("D44").Value = "'5.32"
$ws.Range("E44").Value = "  -4.38%  "

$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("E46").Value = "  +0.32%  "

$ws.Range("D47").Value = "'100.30"
$ws.Range("E47").Value = "  -2.48%  "

$ws.Range("E48").Value = "  +2.67%  "

$ws.Range("E49").Value = "  +1.15%  "

$ws.Range("D50").Value = "'0.438"
$ws.Range("E50").Value = "  -1.24%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.446.02"
$ws.Range("E51").Value = "  +0.36%  "
